$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two SKU columns (E = "TPxxxxx-1.jpg", F = "TPxxxxx"); this shifts
# the image-url / category columns left by two (G->E, H->F, I->G, J->H, K->I, L->J).
$ws.Range("E:F").EntireColumn.Delete()

# New column K holds a literal quantity/invoice value of 100 for every data row.
for ($r = 1; $r -le 19; $r++) {
    $ws.Range("K$r").Value = 100
}

# Reset the sheet selection to a single cell (J1) instead of the whole header row.
$ws.Range("J1").Select()
